# Daily attendance processing - 2026-01-29 06:22:06
# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (column G) wherever both values are present,
# changing "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
